# Insert a new weekly price record at row 38 ("Hortaliza, Comercializadora
# del Agro de Limarí - Haba"). Inserting a whole row shifts the existing
# rows 38-58 down to 39-59, preserving all of their data/formatting, and
# extends the used range to A1:R59 - matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 38..58 down to 39..59.
$ws.Rows.Item(38).Insert()

# Populate the newly inserted row 38 with the new record.
$ws.Range("A38").Value = 2
$ws.Range("B38").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C38").Value = "Coquimbo"
$ws.Range("D38").Value = 44762
$ws.Range("E38").Value = 4
$ws.Range("F38").Value = 100112026
$ws.Range("G38").Value = "Haba"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 1500
$ws.Range("K38").Value = 11000
$ws.Range("L38").Value = 12000
$ws.Range("M38").Value = 11500
$ws.Range("N38").Value = "$/saco 25 kilos"
$ws.Range("O38").Value = "Provincia de Limarí"
$ws.Range("P38").Value = 460
$ws.Range("Q38").Value = 25
$ws.Range("R38").Value = "Hortaliza"
